$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Sheet, [string]$Ref, [string]$NewValue)
    $cell = $Sheet.Range($Ref)
    # Prefix with apostrophe so Excel stores numeric-looking strings as text,
    # then reset the style so no quotePrefix/number-format style sticks to the cell.
    $cell.Value = "'" + $NewValue
    $cell.Style = "Normal"
}

Set-TextValue $ws "D2" "84.187.50"
Set-TextValue $ws "E2" "  +5.58%  "
Set-TextValue $ws "D3" "3.277.04"
Set-TextValue $ws "E3" "  +1.88%  "
Set-TextValue $ws "E4" "  +0.17%  "
Set-TextValue $ws "D5" "219.26"
Set-TextValue $ws "E5" "  +4.14%  "
Set-TextValue $ws "D6" "630.92"
Set-TextValue $ws "E6" "  -1.42%  "
Set-TextValue $ws "D7" "0.318"
Set-TextValue $ws "E7" "  +22.59%  "
Set-TextValue $ws "D8" "0.999"
Set-TextValue $ws "E8" "  -0.01%  "
Set-TextValue $ws "D9" "0.591"
Set-TextValue $ws "E9" "  -1.96%  "
Set-TextValue $ws "D10" "3.275.68"
Set-TextValue $ws "E10" "  +1.92%  "
Set-TextValue $ws "D11" "0.596"
Set-TextValue $ws "E11" "  -1.43%  "
Set-TextValue $ws "D12" "0.0000277"
Set-TextValue $ws "E12" "  +5.27%  "
Set-TextValue $ws "E13" "  -0.08%  "
Set-TextValue $ws "D14" "3.890.09"
Set-TextValue $ws "E14" "  +2.39%  "
Set-TextValue $ws "D15" "5.40"
Set-TextValue $ws "E15" "  -0.62%  "
Set-TextValue $ws "D16" "33.00"
Set-TextValue $ws "E16" "  +1.55%  "
Set-TextValue $ws "D17" "84.304.76"
Set-TextValue $ws "E17" "  +6.03%  "
Set-TextValue $ws "D18" "3.269.92"
Set-TextValue $ws "E18" "  +2.14%  "
Set-TextValue $ws "D19" "3.19"
Set-TextValue $ws "E19" "  +5.56%  "
Set-TextValue $ws "D20" "14.42"
Set-TextValue $ws "E20" "  -1.81%  "
Set-TextValue $ws "D21" "450.39"
Set-TextValue $ws "E21" "  +0.86%  "
Set-TextValue $ws "D22" "9.11"
Set-TextValue $ws "E22" "  -2.84%  "
Set-TextValue $ws "D23" "5.21"
Set-TextValue $ws "E23" "  -1.11%  "
Set-TextValue $ws "D24" "7.47"
Set-TextValue $ws "E24" "  +7.57%  "
Set-TextValue $ws "D25" "5.24"
Set-TextValue $ws "E25" "  +7.92%  "
Set-TextValue $ws "D26" "12.13"
Set-TextValue $ws "E26" "  +10.96%  "
Set-TextValue $ws "D27" "3.436.88"
Set-TextValue $ws "E27" "  +2.00%  "
Set-TextValue $ws "D28" "77.68"
Set-TextValue $ws "E28" "  +0.15%  "
Set-TextValue $ws "E29" "  +0.15%  "
Set-TextValue $ws "E30" "  +0.53%  "
Set-TextValue $ws "D31" "9.20"
Set-TextValue $ws "E31" "  -0.28%  "
Set-TextValue $ws "E32" "  -0.26%  "
Set-TextValue $ws "D33" "0.158"
Set-TextValue $ws "E33" "  +29.70%  "
Set-TextValue $ws "D34" "570.11"
Set-TextValue $ws "E34" "  +1.58%  "
Set-TextValue $ws "E36" "  -3.15%  "
Set-TextValue $ws "D37" "2.01"
Set-TextValue $ws "E37" "  -1.91%  "
Set-TextValue $ws "D38" "23.27"
Set-TextValue $ws "E38" "  +0.23%  "
Set-TextValue $ws "B39" "RenderToken"
Set-TextValue $ws "C39" "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
Set-TextValue $ws "D39" "6.19"
Set-TextValue $ws "E39" "  +8.11%  "
Set-TextValue $ws "B40" "FirstDigitalUSD"
Set-TextValue $ws "C40" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws "D40" "0.998"
Set-TextValue $ws "E40" "  -0.18%  "
Set-TextValue $ws "D41" "0.410"
Set-TextValue $ws "E41" "  -0.98%  "
Set-TextValue $ws "D42" "2.06"
Set-TextValue $ws "E42" "  +12.23%  "
Set-TextValue $ws "B43" "WhiteBITCoin"
Set-TextValue $ws "C43" "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue $ws "D43" "20.93"
Set-TextValue $ws "E43" "  +3.22%  "
Set-TextValue $ws "B44" "dogwifhat"
Set-TextValue $ws "C44" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws "D44" "3.04"
Set-TextValue $ws "E44" "  +12.69%  "
Set-TextValue $ws "D45" "159.63"
Set-TextValue $ws "E45" "  -2.03%  "
Set-TextValue $ws "E46" "  +0.06%  "
Set-TextValue $ws "D47" "190.50"
Set-TextValue $ws "E47" "  -1.98%  "
Set-TextValue $ws "D48" "45.00"
Set-TextValue $ws "E48" "  +4.69%  "
Set-TextValue $ws "D49" "1.33"
Set-TextValue $ws "E49" "  -1.24%  "
Set-TextValue $ws "D50" "0.777"
Set-TextValue $ws "E50" "  -3.20%  "
Set-TextValue $ws "D51" "26.19"
Set-TextValue $ws "E51" "  +0.55%  "
